# database/ezoom_db_objects.xlsx
# Commit: "fixed bugs in new_order_demand_prediction and updated
#          order_dispatch_display_per_product.sql and updated
#          load_procedures and ezoom_db_objects.xlsx"
#
# This adds one new catalogued DB object - a stored procedure named
# order_dispatch_display_per_product.sql - as a new row right before the
# existing "view" rows (which get pushed down by one row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand new row at 125; the existing rows 125-127 (the "view"
# entries) shift down to 126-128.
$ws.Rows.Item(125).Insert()

# Populate the newly inserted row with the new stored procedure entry.
$ws.Range("A125").Value = "stored procedure"
$ws.Range("B125").Value = "order_dispatch_display_per_product.sql"
$ws.Range("C125").Value = "dispatch product per order and predict max final product allowed by inventory"

# Mirror the author's final cursor/viewport position after the edit.
$ws.Range("C131").Select()
